$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the outlier row (6, 0, 3) - row 47 - which shifts subsequent rows up by one
$ws.Rows(47).Delete()

# Update RMSE column (B) values: values now represent the maximum IC-capacity
# per month rather than the sum of IC-capacity
$newValues = @{
    2 = 340.6502365674644
    3 = 380.1709986241312
    4 = 466.4927275488926
    5 = 442.3999701996131
    6 = 282.3504441540706
    7 = 495.8227692109853
    8 = 503.5020521473747
    9 = 534.2329043762988
    10 = 508.2810796541641
    11 = 437.3966666263567
    12 = 337.4029434358142
    13 = 305.4142891313475
    14 = 284.42040542354
    15 = 502.5108446169883
    16 = 349.3221794436784
    17 = 216.557585124399
    18 = 225.4852715383226
    19 = 368.253494630708
    20 = 224.6647038751593
    21 = 309.9864289648673
    22 = 307.4330988189717
    23 = 352.8185281751925
    24 = 304.0039148229766
    25 = 201.0344671326881
    26 = 361.3978834513393
    27 = 265.2880644961264
    28 = 282.6231269891917
    29 = 252.5436773803101
    30 = 334.024038593734
    31 = 349.837404554086
    32 = 307.9019510352074
    33 = 189.6693852564502
    34 = 365.4000240489517
    35 = 287.6846542207747
    36 = 201.3461906860068
    37 = 438.1491704926532
    38 = 439.0007528979557
    39 = 393.599403680087
    40 = 311.2395037424493
    41 = 222.0121201869233
    42 = 380.9072746760096
    43 = 290.6370800647912
    44 = 365.991398300569
    45 = 399.9812871821335
    46 = 566.9913101431812
    47 = 364.2436063050155
    48 = 229.4288501838379
    49 = 248.4104850046107
    50 = 371.2264591495643
    51 = 348.7460421366534
    52 = 344.818675737992
    53 = 385.7467845371036
    54 = 222.9846743578483
    55 = 281.0098902576678
    56 = 242.0748495274166
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
